# Daily attendance processing - 2025-11-09 15:20:15
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: reorder "Recorded By" e-mail list (G2)
$ws.Range("G2").Value = "Veronia.rafat@med.asu.edu.eg, System, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

# Row 3: add recorder and update attendance count (G3, H3)
$ws.Range("G3").Value = "majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("H3").Value = "38/251"

# Row 9: reorder "Recorded By" e-mail list (G9)
$ws.Range("G9").Value = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

# Row 10: updated average attendance % (L10) - force text format so the
# "24.6%" literal is stored verbatim instead of being auto-converted to a
# percentage number, then restore the general number format so the cell's
# style does not otherwise drift from its original appearance.
$cellL10 = $ws.Range("L10")
$cellL10.NumberFormat = "@"
$cellL10.Value = "24.6%"
$cellL10.NumberFormat = "general"

# Row 15: updated average attendance % (S15), same text-literal handling
$cellS15 = $ws.Range("S15")
$cellS15.NumberFormat = "@"
$cellS15.Value = "24.6%"
$cellS15.NumberFormat = "general"

# Row 28: reorder "Recorded By" e-mail list (G28)
$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
